$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.514.64"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "'2.513.98"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'316.92"
$ws.Range("E5").Value = "  +4.84%  "
$ws.Range("D6").Value = "'94.07"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").Value = "'0.575"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").Value = "'2.903.30"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "'2.537.28"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "'15.18"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "'42.646.83"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "'6.66"
$ws.Range("E20").Value = "  +4.43%  "
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "'69.24"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").Value = "'2.96"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "'26.70"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("D29").Value = "'40.65"
$ws.Range("E29").Value = "  +9.17%  "
$ws.Range("D30").Value = "'10.23"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").Value = "'5.94"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "'158.40"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").Value = "'19.15"
$ws.Range("E33").Value = "  +4.28%  "
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("D35").Value = "'3.26"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").Value = "'0.0779"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "'23.37"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("D41").Value = "'2.30"
$ws.Range("E41").Value = "  +17.48%  "
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "'0.0303"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.31"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'3.76"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "'2.011.42"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("D47").Value = "'85.32"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("D48").Value = "'8.86"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'74.37"
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "'2.756.63"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'101.92"
$ws.Range("E51").Value = "  +1.55%  "
